# Update the "Quarterly Data Levels" sheet with revised Household Deposits (K),
# Savings Rate (M) and one Payroll employment (O) figure per the source data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly Data Levels")

$ws.Range("K34").Value = 117353
$ws.Range("K35").Value = 119112
$ws.Range("K36").Value = 125150
$ws.Range("K37").Value = 128661
$ws.Range("K38").Value = 128178
$ws.Range("K39").Value = 131960
$ws.Range("K40").Value = 139025
$ws.Range("K41").Value = 140971
$ws.Range("K42").Value = 145928
$ws.Range("K43").Value = 150228
$ws.Range("K44").Value = 158682
$ws.Range("K45").Value = 165860
$ws.Range("K46").Value = 171719
$ws.Range("K47").Value = 178870
$ws.Range("K48").Value = 184801
$ws.Range("K49").Value = 193464
$ws.Range("K50").Value = 191262
$ws.Range("K51").Value = 189440
$ws.Range("K52").Value = 194190
$ws.Range("K53").Value = 193650
$ws.Range("K54").Value = 191609
$ws.Range("K55").Value = 193711
$ws.Range("K56").Value = 195316
$ws.Range("K57").Value = 197292
$ws.Range("K58").Value = 199389
$ws.Range("K59").Value = 201356
$ws.Range("K60").Value = 205269
$ws.Range("K61").Value = 206864
$ws.Range("K62").Value = 203563
$ws.Range("K63").Value = 198573
$ws.Range("K64").Value = 197020
$ws.Range("K65").Value = 191550
$ws.Range("K66").Value = 185644
$ws.Range("K67").Value = 183608
$ws.Range("K68").Value = 183770
$ws.Range("K69").Value = 185159
$ws.Range("K70").Value = 182862
$ws.Range("K71").Value = 185642
$ws.Range("K72").Value = 187074
$ws.Range("K73").Value = 188273
$ws.Range("K74").Value = 186715
$ws.Range("K75").Value = 183580
$ws.Range("K76").Value = 182392
$ws.Range("K77").Value = 182338
$ws.Range("K78").Value = 181138
$ws.Range("K79").Value = 180972
$ws.Range("K80").Value = 181338
$ws.Range("K81").Value = 183774
$ws.Range("K82").Value = 180251
$ws.Range("K83").Value = 179964
$ws.Range("K84").Value = 182578
$ws.Range("M84").Value = 11.9
$ws.Range("K85").Value = 180782
$ws.Range("K86").Value = 178509
$ws.Range("K87").Value = 179155
$ws.Range("K88").Value = 176672
$ws.Range("K89").Value = 176242
$ws.Range("K90").Value = 174952
$ws.Range("K91").Value = 175585
$ws.Range("M91").Value = 14.2
$ws.Range("K92").Value = 179171
$ws.Range("M92").Value = 14.7
$ws.Range("K93").Value = 184755
$ws.Range("K94").Value = 182457
$ws.Range("M94").Value = 11.5
$ws.Range("K95").Value = 183487
$ws.Range("M95").Value = 12.5
$ws.Range("K96").Value = 188255
$ws.Range("M96").Value = 13.4
$ws.Range("K97").Value = 189240
$ws.Range("K98").Value = 184672
$ws.Range("K99").Value = 188551
$ws.Range("M99").Value = 13.3
$ws.Range("K100").Value = 194972
$ws.Range("M100").Value = 13.8
$ws.Range("K101").Value = 201853
$ws.Range("K102").Value = 200288
$ws.Range("M102").Value = 24.6
$ws.Range("K103").Value = 204997
$ws.Range("K104").Value = 212254
$ws.Range("M104").Value = 21.2
$ws.Range("K105").Value = 214953
$ws.Range("K106").Value = 199768
$ws.Range("M106").Value = 30.9
$ws.Range("K107").Value = 205745
$ws.Range("K108").Value = 208146
$ws.Range("M108").Value = 19.7
$ws.Range("K109").Value = 209529
$ws.Range("M109").Value = 19.3
$ws.Range("K110").Value = 207735
$ws.Range("M110").Value = 16.3
$ws.Range("K111").Value = 211523
$ws.Range("M111").Value = 16.3
$ws.Range("K112").Value = 210920
$ws.Range("M112").Value = 15.1
$ws.Range("K113").Value = 218475
$ws.Range("M113").Value = 14.1
$ws.Range("K114").Value = 218651
$ws.Range("M114").Value = 15.3
$ws.Range("K115").Value = 216603
$ws.Range("M115").Value = 13.8
$ws.Range("K116").Value = 218999
$ws.Range("M116").Value = 13.4
$ws.Range("K117").Value = 227916
$ws.Range("M117").Value = 13
$ws.Range("K118").Value = 227023
$ws.Range("M118").Value = 13.9
$ws.Range("K119").Value = 233944
$ws.Range("M119").Value = 12.9
$ws.Range("K120").Value = 242926
$ws.Range("M120").Value = 14.1
$ws.Range("O120").Value = 2505933.33333333
